$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C4").Value = "https://logs.xiti.com/hit.xiti?s=388889&s2=1&p=particuliers::accueil_PARTICULIERS&hl=23x52x18&xtor=sec-322800100-goo-[retraite_generique%20%2B%20marque_bmm]-[103507913820]-s-[%2Bperp%20%2Bla%20%2Bposte]&roinbh=392&ac=&an=&lng=en-US&idp=2352183188828&jv=0&re=625x803&vtag=44007&hl=23x52x18&r=1440x900x24x24&ref="
$ws.Range("C5").Value = "https://logs.xiti.com/hit.xiti?s=388889&s2=21&p=Outils::aide::cookies&hl=23x52x56&xtor=sec-322800100-goo-[retraite_generique%20%2B%20marque_bmm]-[103507913820]-s-[%2Bperp%20%2Bla%20%2Bposte]&roinbh=392&ac=&an=&lng=en-US&idp=2352563949647&jv=0&re=625x803&vtag=44007&hl=23x52x56&r=1440x900x24x24&ref="
$ws.Range("C6").Value = "https://logs.xiti.com/hit.xiti?s=388889&s2=21&p=projets&hl=23x53x12&xtor=sec-322800100-goo-[retraite_generique%20%2B%20marque_bmm]-[103507913820]-s-[%2Bperp%20%2Bla%20%2Bposte]&roinbh=392&ac=&an=&lng=en-US&idp=2353124427578&jv=0&re=625x803&vtag=44007&hl=23x53x12&r=1440x900x24x24&ref="

$ws.Range("C2").Copy()
$ws.Range("C4:C6").PasteSpecial(-4122)

$ws.Range("G29").Select()
